## 📊 Actualización automática del dashboard
## Inserta el nuevo comentario de TikTok (cid 7601898749665395464) en la fila 11
## de "Comentarios" y actualiza los totales dependientes en las demás hojas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentarios")

# Inserta una fila nueva encima de la fila 11 (desplaza el resto hacia abajo)
$ws.Rows.Item(11).Insert()

# --- Rellena la fila 11 con el nuevo comentario ---
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "TikTok"
$ws.Range("C11").Value = "https://vt.tiktok.com/ZSannqcDU/"
$ws.Range("D11").Value = "https://vt.tiktok.com/ZSannqcDU/"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "Busca en Youtube, escribe aqua mochis y te aparece"
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = ""

$ws.Range("I11").Value = 0

# comment_id se guarda como texto (el id es demasiado largo para numero sin perder precision).
# Truco: escribir una formula de texto en una celda lejos de los datos, copiarla y pegar solo
# los valores en J11 -- asi el tipo queda como texto sin tener que tocar el NumberFormat
# de la celda (lo que evitaria crear un estilo nuevo de mas).
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="7601898749665395464"'
$scratch.Copy()
$ws.Range("J11").PasteSpecial(-4163)
$scratch.EntireColumn.Delete()

$ws.Range("K11").Value = 1769954987

$ws.Range("L11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L11").Value = 46054.59012731481

$ws.Range("M11").NumberFormat = "yyyy-mm-dd"
$ws.Range("M11").NumberFormat = "YYYY-MM-DD"
$ws.Range("M11").Value = 46054

$ws.Range("N11").Value = "14:09:47"
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = ""
$ws.Range("Q11").Value = "https://www.tiktok.com/@"
$ws.Range("R11").Value = ""
$ws.Range("S11").Value = "{'videoWebUrl': 'https://www.tiktok.com/@alpinacol/video/7600955886223232277', 'submittedVideoUrl': 'https://vt.tiktok.com/ZSannqcDU/', 'input': 'https://vt.tiktok.com/ZSannqcDU/', 'cid': '7601898749665395464', 'createTime': 1769954987, 'createTimeISO': '2026-02-01T14:09:47.000Z', 'text': 'Busca en Youtube, escribe aqua mochis y te aparece', 'diggCount': 0, 'likedByAuthor': False, 'repliesToId': '7601130474090464007', 'replyCommentTotal': None, 'uid': '7462023650225783816', 'uniqueId': 'salome_."

# --- Actualiza los totales en las demas hojas (41 -> 42 comentarios) ---
$resumen = $wb.Worksheets.Item("Resumen_Posts")
$resumen.Range("D2").Value = 42

$stats = $wb.Worksheets.Item("Stats_Plataforma")
$stats.Range("C2").Value = 42
$stats.Range("D2").Value = 0.71

$replies = $wb.Worksheets.Item("Stats_Replies")
$replies.Range("B2").Value = 42

Write-Output "dashboard actualizado"
